$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set columns D and E for changed rows to Text format first so that
# numeric-looking / percent-looking strings are preserved verbatim as text,
# matching the source data (inline strings), not converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.627.80'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +7.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.738.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9952'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +5.83%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3724'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.70'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +5.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3389'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.207'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07539'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9929'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.440'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +6.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.57'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.021'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.734.89'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001100'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06677'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.62'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9935'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.82'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.184'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.15'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '26.568.74'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.465'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.42%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +5.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.428'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +16.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.01'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.66'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.920.44'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '131.89'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.118'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.230'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08594'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.25'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +8.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.701'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.475'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06371'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02360'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2177'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.696'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.245'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6272'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +5.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.55'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +13.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9940'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.905'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6080'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.64'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.071'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07344'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.14'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.37%  '
